# CORDEX_simulations_ATLAS.xlsx - "Index numbers updated for all the domains"
#
# This script renumbers the index column (A) on several domain sheets,
# clears the leftover highlight formatting on a handful of cells whose
# index value stopped being "new" (so they revert to the plain/unstyled
# look used by the rest of the column), and replays the cursor/selection
# state that was left behind in each sheet, finishing with "14.SEA" as the
# active tab (matching the saved workbookView.activeTab=10).

$wb = $excel.ActiveWorkbook

function Clear-IndexHighlight($ws, $targetCell, $sourceCell) {
    # Copies ONLY the formatting (not the value) from a plain, already
    # "s=1"-styled cell onto $targetCell, so the renumbered cell matches
    # the look of the rest of the index column.
    $ws.Range($sourceCell).Copy() | Out-Null
    $ws.Range($targetCell).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# 1.SAM-final — selection only
# ---------------------------------------------------------------------
$wsSAM = $wb.Worksheets.Item("1.SAM-final")
$wsSAM.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 2.CAM-final — renumber A14:A26, selection
# ---------------------------------------------------------------------
$wsCAM = $wb.Worksheets.Item("2.CAM-final")
$wsCAM.Range("A14").Value = 9
$wsCAM.Range("A15").Value = 10
$wsCAM.Range("A16").Value = 11
$wsCAM.Range("A17").Value = 12
$wsCAM.Range("A18").Value = 13
$wsCAM.Range("A19").Value = 14
$wsCAM.Range("A20").Value = 15
$wsCAM.Range("A21").Value = 16
$wsCAM.Range("A22").Value = 17
$wsCAM.Range("A23").Value = 18
$wsCAM.Range("A24").Value = 19
$wsCAM.Range("A25").Value = 20
$wsCAM.Range("A26").Value = 21
$wsCAM.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 3.NAM-final — renumber A16:A25, selection
# ---------------------------------------------------------------------
$wsNAM = $wb.Worksheets.Item("3.NAM-final")
$wsNAM.Range("A16").Value = 11
$wsNAM.Range("A17").Value = 12
$wsNAM.Range("A18").Value = 13
$wsNAM.Range("A19").Value = 14
$wsNAM.Range("A20").Value = 15
$wsNAM.Range("A21").Value = 16
$wsNAM.Range("A22").Value = 17
$wsNAM.Range("A23").Value = 18
$wsNAM.Range("A24").Value = 19
$wsNAM.Range("A25").Value = 20
$wsNAM.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 4.EUR — scroll/selection, and it stops being the active tab
# ---------------------------------------------------------------------
$wsEUR = $wb.Worksheets.Item("4.EUR")
$wsEUR.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 5.AFR — renumber A27:A42, clear stale highlight on A12:A14/A24:A25/A40,
# selection
# ---------------------------------------------------------------------
$wsAFR = $wb.Worksheets.Item("5.AFR")
Clear-IndexHighlight $wsAFR "A12" "A11"
Clear-IndexHighlight $wsAFR "A13" "A11"
Clear-IndexHighlight $wsAFR "A14" "A11"
Clear-IndexHighlight $wsAFR "A24" "A23"
Clear-IndexHighlight $wsAFR "A25" "A23"
$wsAFR.Range("A27").Value = 22
$wsAFR.Range("A28").Value = 23
$wsAFR.Range("A29").Value = 24
$wsAFR.Range("A30").Value = 25
$wsAFR.Range("A31").Value = 26
$wsAFR.Range("A32").Value = 27
$wsAFR.Range("A33").Value = 28
$wsAFR.Range("A34").Value = 29
$wsAFR.Range("A35").Value = 30
$wsAFR.Range("A36").Value = 31
$wsAFR.Range("A37").Value = 32
$wsAFR.Range("A38").Value = 33
$wsAFR.Range("A39").Value = 34
Clear-IndexHighlight $wsAFR "A40" "A39"
$wsAFR.Range("A40").Value = 35
$wsAFR.Range("A41").Value = 36
$wsAFR.Range("A42").Value = 37
$wsAFR.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 6.WAS-final — selection only
# ---------------------------------------------------------------------
$wsWAS = $wb.Worksheets.Item("6.WAS-final")
$wsWAS.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 7.EAS-final — selection only
# ---------------------------------------------------------------------
$wsEAS = $wb.Worksheets.Item("7.EAS-final")
$wsEAS.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 9.AUS — selection only
# ---------------------------------------------------------------------
$wsAUS = $wb.Worksheets.Item("9.AUS")
$wsAUS.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 10.ANT — selection only
# ---------------------------------------------------------------------
$wsANT = $wb.Worksheets.Item("10.ANT")
$wsANT.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 11.ARC-final — selection only
# ---------------------------------------------------------------------
$wsARC = $wb.Worksheets.Item("11.ARC-final")
$wsARC.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 14.SEA — clear stale highlight on A12, selection, and becomes the
# active tab (workbookView.activeTab=10 / tabSelected=true)
# ---------------------------------------------------------------------
$wsSEA = $wb.Worksheets.Item("14.SEA")
Clear-IndexHighlight $wsSEA "A12" "A11"
$wsSEA.Range("B20").Select() | Out-Null
$wsSEA.Activate() | Out-Null
